# Edit: insert a new data row into the weekly Hortaliza (Repollo) price table.
# The new row is inserted at row 1082, pushing the existing rows 1082-1197
# down to 1083-1198 (dimension grows from A1:R1197 to A1:R1198).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 1082 (existing row 1082 and below shift down by one).
$ws.Rows.Item(1082).Insert()

# Populate the newly inserted row 1082 with the new record.
$ws.Cells.Item(1082, 1).Value = 10
$ws.Cells.Item(1082, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1082, 3).Value = "La Araucanía"
$ws.Cells.Item(1082, 4).Value = 45212
$ws.Cells.Item(1082, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1082, 5).Value = 9
$ws.Cells.Item(1082, 6).Value = 100112006
$ws.Cells.Item(1082, 7).Value = "Repollo"
$ws.Cells.Item(1082, 8).Value = "Crespo record"
$ws.Cells.Item(1082, 9).Value = "Primera"
$ws.Cells.Item(1082, 10).Value = 600
$ws.Cells.Item(1082, 11).Value = 1000
$ws.Cells.Item(1082, 12).Value = 1000
$ws.Cells.Item(1082, 13).Value = 1000
$ws.Cells.Item(1082, 14).Value = "`$/unidad"
$ws.Cells.Item(1082, 15).Value = "Región del Maule"
$ws.Cells.Item(1082, 16).Value = 1000
$ws.Cells.Item(1082, 17).Value = 1
$ws.Cells.Item(1082, 18).Value = "Hortaliza"
